# Update countries & provincias Spain
#
# This updates the "Pais" (country) COVID stats sheet with a newer data
# pull:
#   - Brasil (row 11)   gets refreshed totals
#   - Sri Lanka (row 105) gets refreshed totals
#   - Honduras overtakes Bulgaria in total cases, so the two rows that used
#     to read Bulgaria/Honduras (rows 79/80, sorted by total cases
#     descending) now read Honduras/Bulgaria; Honduras' row gets the new
#     numbers while Bulgaria keeps its previous (unchanged) figures.
#   - The "datos actualizados" timestamp footer is bumped forward half an
#     hour, from 03:35 to 04:05.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: Brasil ---------------------------------------------------
$ws.Range("B11").Value = 169594
$ws.Range("C11").Value = 451
$ws.Range("E11").Value = 90557
$ws.Range("G11").Value = 28
$ws.Range("H11").Value = 11653

# --- Rows 79/80: Honduras overtakes Bulgaria ---------------------------
# Row 79 now holds Honduras (previously Bulgaria) with fresh numbers.
$ws.Range("A79").Value = "Honduras"
$ws.Range("B79").Value = 2100
$ws.Range("C79").Value = 128
$ws.Range("D79").Value = 206
$ws.Range("E79").Value = 1778
$ws.Range("F79").Value = 10
$ws.Range("G79").Value = 8
$ws.Range("H79").Value = 116

# Row 80 now holds Bulgaria (previously Honduras) with its old, unchanged
# figures (Bulgaria's data did not move this update).
$ws.Range("A80").Value = "Bulgaria"
$ws.Range("B80").Value = 1990
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 461
$ws.Range("E80").Value = 1436
$ws.Range("F80").Value = 50
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 93

# --- Row 105: Sri Lanka -------------------------------------------------
$ws.Range("B105").Value = 869
$ws.Range("C105").Value = 6
$ws.Range("E105").Value = 517

# --- Timestamp footer -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 04:05"
